$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: De'Anthony Melton
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 'De''Anthony Melton'
$ws.Range("D2").Value = 'PG'
$ws.Range("E2").Value = '6-2'
$ws.Range("F2").Value = 200
$ws.Range("G2").Value = 'May 28, 1998'
$ws.Range("H2").Value = 'us'
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 'USC'
$ws.Range("K2").Value = 'https://www.basketball-reference.com/players/m/meltode01.html'

# Row 3: Georges Niang
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 'Georges Niang'
$ws.Range("D3").Value = 'PF'
$ws.Range("E3").Value = '6-7'
$ws.Range("F3").Value = 230
$ws.Range("G3").Value = 'June 17, 1993'
$ws.Range("H3").Value = 'us'
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 'Iowa State'
$ws.Range("K3").Value = 'https://www.basketball-reference.com/players/n/niangge01.html'

# Row 4: P.J. Tucker
$ws.Range("B4").Value = 17
$ws.Range("C4").Value = 'P.J. Tucker'
$ws.Range("D4").Value = 'PF'
$ws.Range("E4").Value = '6-5'
$ws.Range("F4").Value = 245
$ws.Range("G4").Value = 'May 5, 1985'
$ws.Range("H4").Value = 'us'
$ws.Range("I4").Value = 11
$ws.Range("J4").Value = 'Texas'
$ws.Range("K4").Value = 'https://www.basketball-reference.com/players/t/tuckepj01.html'

# Row 5: Shake Milton
$ws.Range("B5").Value = 18
$ws.Range("C5").Value = 'Shake Milton'
$ws.Range("D5").Value = 'SG'
$ws.Range("E5").Value = '6-5'
$ws.Range("F5").Value = 205
$ws.Range("G5").Value = 'September 26, 1996'
$ws.Range("H5").Value = 'us'
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 'SMU'
$ws.Range("K5").Value = 'https://www.basketball-reference.com/players/m/miltosh01.html'

# Row 6: Tobias Harris
$ws.Range("B6").Value = 12
$ws.Range("C6").Value = 'Tobias Harris'
$ws.Range("D6").Value = 'PF'
$ws.Range("E6").Value = '6-8'
$ws.Range("F6").Value = 226
$ws.Range("G6").Value = 'July 15, 1992'
$ws.Range("H6").Value = 'us'
$ws.Range("I6").Value = 11
$ws.Range("J6").Value = 'Tennessee'
$ws.Range("K6").Value = 'https://www.basketball-reference.com/players/h/harrito02.html'

# Row 7: Joel Embiid
$ws.Range("B7").Value = 21
$ws.Range("C7").Value = 'Joel Embiid'
$ws.Range("D7").Value = 'C'
$ws.Range("E7").Value = '7-0'
$ws.Range("F7").Value = 280
$ws.Range("G7").Value = 'March 16, 1994'
$ws.Range("H7").Value = 'cm'
$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 'Kansas'
$ws.Range("K7").Value = 'https://www.basketball-reference.com/players/e/embiijo01.html'

# Row 8: Paul Reed
$ws.Range("B8").Value = 44
$ws.Range("C8").Value = 'Paul Reed'
$ws.Range("D8").Value = 'C'
$ws.Range("E8").Value = '6-9'
$ws.Range("F8").Value = 210
$ws.Range("G8").Value = 'June 14, 1999'
$ws.Range("H8").Value = 'us'
$ws.Range("I8").Value = 2
$ws.Range("J8").Value = 'DePaul'
$ws.Range("K8").Value = 'https://www.basketball-reference.com/players/r/reedpa01.html'

# Row 9: Montrezl Harrell
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 'Montrezl Harrell'
$ws.Range("D9").Value = 'C'
$ws.Range("E9").Value = '6-7'
$ws.Range("F9").Value = 240
$ws.Range("G9").Value = 'January 26, 1994'
$ws.Range("H9").Value = 'us'
$ws.Range("I9").Value = 7
$ws.Range("J9").Value = 'Louisville'
$ws.Range("K9").Value = 'https://www.basketball-reference.com/players/h/harremo01.html'

# Row 10: James Harden
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 'James Harden'
$ws.Range("D10").Value = 'PG'
$ws.Range("E10").Value = '6-5'
$ws.Range("F10").Value = 220
$ws.Range("G10").Value = 'August 26, 1989'
$ws.Range("H10").Value = 'us'
$ws.Range("I10").Value = 13
$ws.Range("J10").Value = 'Arizona State'
$ws.Range("K10").Value = 'https://www.basketball-reference.com/players/h/hardeja01.html'

# Row 11: Tyrese Maxey
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 'Tyrese Maxey'
$ws.Range("D11").Value = 'SG'
$ws.Range("E11").Value = '6-2'
$ws.Range("F11").Value = 200
$ws.Range("G11").Value = 'November 4, 2000'
$ws.Range("H11").Value = 'us'
$ws.Range("I11").Value = 2
$ws.Range("J11").Value = 'Kentucky'
$ws.Range("K11").Value = 'https://www.basketball-reference.com/players/m/maxeyty01.html'

# Row 12: Danuel House Jr.
$ws.Range("B12").Value = 25
$ws.Range("C12").Value = 'Danuel House Jr.'
$ws.Range("D12").Value = 'SF'
$ws.Range("E12").Value = '6-6'
$ws.Range("F12").Value = 220
$ws.Range("G12").Value = 'June 7, 1993'
$ws.Range("H12").Value = 'us'
$ws.Range("I12").Value = 6
$ws.Range("J12").Value = 'Texas A&M'
$ws.Range("K12").Value = 'https://www.basketball-reference.com/players/h/houseda01.html'

# Row 13: Furkan Korkmaz
$ws.Range("B13").Value = 30
$ws.Range("C13").Value = 'Furkan Korkmaz'
$ws.Range("D13").Value = 'SG'
$ws.Range("E13").Value = '6-7'
$ws.Range("F13").Value = 202
$ws.Range("G13").Value = 'July 24, 1997'
$ws.Range("H13").Value = 'tr'
$ws.Range("I13").Value = 5
$ws.Range("K13").Value = 'https://www.basketball-reference.com/players/k/korkmfu01.html'

# Row 14: Jalen McDaniels
$ws.Range("B14").Value = 7
$ws.Range("C14").Value = 'Jalen McDaniels'
$ws.Range("D14").Value = 'SF'
$ws.Range("E14").Value = '6-9'
$ws.Range("F14").Value = 205
$ws.Range("G14").Value = 'January 31, 1998'
$ws.Range("H14").Value = 'us'
$ws.Range("I14").Value = 3
$ws.Range("J14").Value = 'San Diego State'
$ws.Range("K14").Value = 'https://www.basketball-reference.com/players/m/mcdanja01.html'

# Row 15: Jaden Springer
$ws.Range("B15").Value = 11
$ws.Range("C15").Value = 'Jaden Springer'
$ws.Range("D15").Value = 'SG'
$ws.Range("E15").Value = '6-4'
$ws.Range("F15").Value = 204
$ws.Range("G15").Value = 'September 25, 2002'
$ws.Range("H15").Value = 'us'
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 'Tennessee'
$ws.Range("K15").Value = 'https://www.basketball-reference.com/players/s/sprinja01.html'

# Row 16: Louis King (TW)
$ws.Range("C16").Value = 'Louis King (TW)'
$ws.Range("D16").Value = 'SF'
$ws.Range("E16").Value = '6-7'
$ws.Range("F16").Value = 205
$ws.Range("G16").Value = 'April 6, 1999'
$ws.Range("H16").Value = 'us'
$ws.Range("I16").Value = 3
$ws.Range("J16").Value = 'Oregon'
$ws.Range("K16").Value = 'https://www.basketball-reference.com/players/k/kinglo02.html'

# Row 17: Mac McClung (TW)
$ws.Range("C17").Value = 'Mac McClung (TW)'
$ws.Range("D17").Value = 'SG'
$ws.Range("E17").Value = '6-2'
$ws.Range("F17").Value = 185
$ws.Range("G17").Value = 'January 6, 1999'
$ws.Range("H17").Value = 'us'
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 'Georgetown, Texas Tech'
$ws.Range("K17").Value = 'https://www.basketball-reference.com/players/m/mccluma01.html'

# Row 18: Dewayne Dedmon
$ws.Range("C18").Value = 'Dewayne Dedmon'
$ws.Range("D18").Value = 'C'
$ws.Range("E18").Value = '7-0'
$ws.Range("F18").Value = 245
$ws.Range("G18").Value = 'August 12, 1989'
$ws.Range("H18").Value = 'us'
$ws.Range("I18").Value = 9
$ws.Range("J18").Value = 'USC'
$ws.Range("K18").Value = 'https://www.basketball-reference.com/players/d/dedmode01.html'
